{"js": "// Commit \"Add files via upload\" re-saved the resume with WPS Office, which\n// merges adjacent runs that share formatting and drops the <w:proofErr/>\n// spell-check markers it doesn't use \u2014 none of that changes the document's\n// visible text. The one real content edit buried in the diff is in the\n// summary paragraph: \"Results-driven Data Engineer with 2.7 years of\n// experience\" becomes \"... 2.8 years of experience\". Reproduce that.\n\nconst body = context.document.body;\nconst matches = body.search(\"Data Engineer with 2.7\", { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error('Expected to find \"Data Engineer with 2.7\" in the document body.');\n}\n\nfor (let i = 0; i < matches.items.length; i++) {\n  matches.items[i].insertText(\"Data Engineer with 2.8\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Commit \"Add files via upload\" re-saved the resume with WPS Office, which\n# merges adjacent runs that share formatting and drops the <w:proofErr/>\n# spell-check markers it doesn't use -- none of that changes the document's\n# visible text. The one real content edit buried in the diff is in the\n# summary paragraph: \"Results-driven Data Engineer with 2.7 years of\n# experience\" becomes \"... 2.8 years of experience\". Reproduce that.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Data Engineer with 2.7\"\n$find.Replacement.Text = \"Data Engineer with 2.8\"\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$find.Execute(\n    [ref]$find.Text,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    $wdFindContinue,\n    $false,\n    [ref]$find.Replacement.Text,\n    $wdReplaceAll\n)\n"}
